# Indicateurs2.xlsx - "Add files via upload" re-edit
#
# The sheet's row 2 (the data row) gets several values revised, a new
# row 5 is added with a single value in column A, and the view's
# scroll position / selection are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing figures in row 2 --------------------------------
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 330
$ws.Range("O2").Value = 200
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 400
$ws.Range("T2").Value = 700

# --- Add a new row further down with a single figure ------------------
$ws.Range("A5").Value = 160

# --- Restore/update the view: scroll position and selection -----------
# Scroll so column H is the first visible column (top-left of the pane).
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1

# Select I2:K2, with I2 as the active cell.
$ws.Range("I2:K2").Select()
